$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update period-header labels (si 35/36/37): shift fiscal quarter from 12-month to 03-month year-end
$ws.Range("C2").Value = "2015/03  (IFRS연결)"
$ws.Range("C3").Value = "2016/03  (IFRS연결)"
$ws.Range("C4").Value = "2017/03  (IFRS연결)"

# Update financial figures for rows 2-6 (restated/corrected values)
$ws.Range("D2").Value = 1582
$ws.Range("E2").Value = 13
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 41
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 28
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 6231
$ws.Range("L2").Value = 2502
$ws.Range("M2").Value = 3728
$ws.Range("N2").Value = 3650
$ws.Range("O2").Value = 78
$ws.Range("P2").Value = 541
$ws.Range("Q2").Value = 221
$ws.Range("R2").Value = -144
$ws.Range("S2").Value = -171
$ws.Range("T2").Value = 90
$ws.Range("U2").Value = 131
$ws.Range("V2").Value = 220
$ws.Range("W2").Value = 0.8100000000000001
$ws.Range("X2").Value = 2.16
$ws.Range("Y2").Value = 0.72
$ws.Range("Z2").Value = 0.54
$ws.Range("AA2").Value = 67.11
$ws.Range("AB2").Value = 585.03
$ws.Range("AC2").Value = 256
$ws.Range("AD2").Value = 67.09999999999999
$ws.Range("AE2").Value = 47247
$ws.Range("AF2").Value = 0.36
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 3.49
$ws.Range("AI2").Value = 167.05
$ws.Range("AJ2").Value = 10821611
$ws.Range("D3").Value = 1583
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 7
$ws.Range("G3").Value = 97
$ws.Range("H3").Value = 64
$ws.Range("I3").Value = 55
$ws.Range("J3").Value = 9
$ws.Range("K3").Value = 6260
$ws.Range("L3").Value = 2494
$ws.Range("M3").Value = 3766
$ws.Range("N3").Value = 3677
$ws.Range("O3").Value = 88
$ws.Range("P3").Value = 541
$ws.Range("Q3").Value = 61
$ws.Range("R3").Value = 80
$ws.Range("S3").Value = -132
$ws.Range("T3").Value = 366
$ws.Range("U3").Value = -305
$ws.Range("V3").Value = 427
$ws.Range("W3").Value = 0.41
$ws.Range("X3").Value = 4.02
$ws.Range("Y3").Value = 1.49
$ws.Range("Z3").Value = 1.02
$ws.Range("AA3").Value = 66.23999999999999
$ws.Range("AB3").Value = 588.78
$ws.Range("AC3").Value = 505
$ws.Range("AD3").Value = 27.91
$ws.Range("AE3").Value = 47597
$ws.Range("AF3").Value = 0.3
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 2.84
$ws.Range("AI3").Value = 56.51
$ws.Range("AJ3").Value = 10821611
$ws.Range("D4").Value = 1373
$ws.Range("E4").Value = -84
$ws.Range("F4").Value = -84
$ws.Range("G4").Value = -69
$ws.Range("H4").Value = -70
$ws.Range("I4").Value = -82
$ws.Range("J4").Value = 12
$ws.Range("K4").Value = 6841
$ws.Range("L4").Value = 3159
$ws.Range("M4").Value = 3682
$ws.Range("N4").Value = 3581
$ws.Range("O4").Value = 101
$ws.Range("P4").Value = 541
$ws.Range("Q4").Value = -223
$ws.Range("R4").Value = -406
$ws.Range("S4").Value = 616
$ws.Range("T4").Value = 810
$ws.Range("U4").Value = -1033
$ws.Range("V4").Value = 1288
$ws.Range("W4").Value = -6.14
$ws.Range("X4").Value = -5.11
$ws.Range("Y4").Value = -2.26
$ws.Range("Z4").Value = -1.07
$ws.Range("AA4").Value = 85.8
$ws.Range("AB4").Value = 567.76
$ws.Range("AC4").Value = -756
$ws.Range("AD4").Value = -17.32
$ws.Range("AE4").Value = 46345
$ws.Range("AF4").Value = 0.28
$ws.Range("AG4").Value = 250
$ws.Range("AH4").Value = 1.91
$ws.Range("AI4").Value = -23.59
$ws.Range("AJ4").Value = 10821611
$ws.Range("D5").Value = 947
$ws.Range("E5").Value = -131
$ws.Range("F5").Value = -131
$ws.Range("G5").Value = -150
$ws.Range("H5").Value = -156
$ws.Range("I5").Value = -165
$ws.Range("J5").Value = 9
$ws.Range("K5").Value = 7199
$ws.Range("L5").Value = 3703
$ws.Range("M5").Value = 3497
$ws.Range("N5").Value = 3389
$ws.Range("O5").Value = 108
$ws.Range("P5").Value = 541
$ws.Range("Q5").Value = -141
$ws.Range("R5").Value = -311
$ws.Range("S5").Value = 448
$ws.Range("T5").Value = 294
$ws.Range("U5").Value = -435
$ws.Range("V5").Value = 1642
$ws.Range("W5").Value = -13.78
$ws.Range("X5").Value = -16.48
$ws.Range("Y5").Value = -4.72
$ws.Range("Z5").Value = -2.22
$ws.Range("AA5").Value = 105.89
$ws.Range("AB5").Value = 533.58
$ws.Range("AC5").Value = -1521
$ws.Range("AD5").Value = -7.92
$ws.Range("AE5").Value = 43861
$ws.Range("AF5").Value = 0.27
$ws.Range("AG5").Value = 200
$ws.Range("AH5").Value = 1.66
$ws.Range("AI5").Value = -9.390000000000001
$ws.Range("AJ5").Value = 10821611
$ws.Range("D6").Value = 1200
$ws.Range("E6").Value = -184
$ws.Range("F6").Value = -184
$ws.Range("G6").Value = -302
$ws.Range("H6").Value = -296
$ws.Range("I6").Value = -299
$ws.Range("K6").Value = 7002
$ws.Range("L6").Value = 3903
$ws.Range("M6").Value = 3099
$ws.Range("N6").Value = 2990
$ws.Range("P6").Value = 541
$ws.Range("Q6").Value = -297
$ws.Range("R6").Value = 68
$ws.Range("S6").Value = 222
$ws.Range("T6").Value = 31
$ws.Range("U6").Value = -328
$ws.Range("V6").Value = 1757
$ws.Range("W6").Value = -15.35
$ws.Range("X6").Value = -24.69
$ws.Range("Y6").Value = -9.359999999999999
$ws.Range("Z6").Value = -4.17
$ws.Range("AA6").Value = 125.96
$ws.Range("AB6").Value = 483.81
$ws.Range("AC6").Value = -2758
$ws.Range("AD6").Value = -3.19
$ws.Range("AE6").Value = 38700
$ws.Range("AF6").Value = 0.23
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 1.7
$ws.Range("AI6").Value = -3.88
$ws.Range("AJ6").Value = 10821611

# Rows 7-9: clear stale data, keep only identifier columns (A, B, C)
$ws.Range("D7:AJ9").ClearContents()
